$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("マスタリスト")

# Delete rows 6 through 74 (inclusive), leaving row 75 to become the new row 6.
$ws.Range("A6:A74").EntireRow.Delete()

# Update the print area to match the new data extent.
$ws.PageSetup.PrintArea = "$A$1:$O$5"
